$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # linear_increases
$ws2 = $wb.Worksheets.Item(2)   # Linear_decrease

# ---------------------------------------------------------------------------
# 1. New "start at 100" budget helper column (F9:F14) + data columns C & D
#    on the linear_increases sheet.
# ---------------------------------------------------------------------------

# Text note cells - entered in this order so the shared-string table grows
# in the same sequence as the authored workbook.
$ws1.Range("F11").Value2 = "The equation still holds and produces 1275"
$ws1.Range("F9").Value2  = "If I want the budget to start from 100 (rather than 0) but still add up to 10,000"
$ws1.Range("F10").Value2 = "5000 = x + 2x + 3x" + [char]0x2026 + "50x"
$ws1.Range("F12").Value2 = "But now it is 5000 = 1275x"
$ws1.Range("F13").Value2 = "Therefore x = 5000/1275"

# Constant used by column D (the "start from 100" increment).
$ws1.Range("F14").Formula = "=5000/1275"

# Column C: cumulative budget starting at 200, step = $B$6 (10000/1275).
$ws1.Range("C6").Value2 = 200
$ws1.Range("C7").Formula = "=C6+`$B`$6"
for ($r = 8; $r -le 55; $r++) {
    $ws1.Cells.Item($r, 3).Formula = "=C" + ($r - 1) + "+`$B`$6"
}

# Column D: cumulative budget starting at 100, step = $F$14 (5000/1275).
$ws1.Range("D6").Value2 = 100
$ws1.Range("D7").Formula = "=D6+`$F`$14"
for ($r = 8; $r -le 55; $r++) {
    $ws1.Cells.Item($r, 4).Formula = "=D" + ($r - 1) + "+`$F`$14"
}

# ---------------------------------------------------------------------------
# 2. Move the existing chart down to make room, and add a new chart that
#    plots the new column C series.
# ---------------------------------------------------------------------------

$co1 = $ws1.ChartObjects().Item(1)
$co1.Left = 385.5185546875
$co1.Top = 258.75
$co1.Width = 443.5
$co1.Height = 216

$co2 = $ws1.ChartObjects().Add(839.2185546875, 265.5, 433.0625, 216)
$chart2 = $co2.Chart
$chart2.ChartType = 4
$chart2.SeriesCollection().NewSeries()
$ser2 = $chart2.SeriesCollection().Item(1)
$ser2.Values = $ws1.Range("C6:C55")

# ---------------------------------------------------------------------------
# 3. View state: linear_increases becomes the active / selected tab, with a
#    new selection, while Linear_decrease loses its "tabSelected" flag.
# ---------------------------------------------------------------------------

$ws1.Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("D7").Select()
